$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H ("Codice Applicativo" stays in G,
# "Nome applicativo" / "Nome Collection Postman" shift right by one).
$ws.Columns.Item(8).Insert()

# New header for the inserted column H.
$ws.Range("H1").Value = "Codice Laboratorio OMR"

# New numeric "Codice Laboratorio OMR" values for rows 2-6.
$ws.Range("H2").Value = 5072024
$ws.Range("H3").Value = 5072025
$ws.Range("H4").Value = 5072026
$ws.Range("H5").Value = 5072027
$ws.Range("H6").Value = 5072028

# The former "Nome Collection Postman" column (now J) gets renamed
# from "Postman Hospital N" to "Postman Hospital Lab N".
$ws.Range("J2").Value = "Postman Hospital Lab 1 "
$ws.Range("J3").Value = "Postman Hospital Lab 2"
$ws.Range("J4").Value = "Postman Hospital Lab 3"
$ws.Range("J5").Value = "Postman Hospital Lab 4"
$ws.Range("J6").Value = "Postman Hospital Lab 5"

# Adjust the view: scroll so column E is the left-most visible column,
# and leave the active selection on J8 (below the data, as in the source).
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("J8").Select()
